$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Pause" value for row 7 (C7), matching the number format (time
# format) used by the other rows in this column.
$ws.Range("C7").NumberFormat = $ws.Range("B7").NumberFormat
$ws.Range("C7").Value = "12:03 - 12:36"

# Update the active selection as recorded after the edit.
$ws.Range("D17").Select()
